$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---
$ws.Name = "Sagarejo"

# --- Row 6 (Urban): years 2010-2022 become confidential/unavailable ("..."),
#     only 2023 (column O) keeps its reported value of 4. ---
$ws.Range("B6:N6").Value = "..."
$ws.Range("O6").Value = 4

# --- Row 7 (Rural): selected years become confidential/unavailable ("..."),
#     the rest keep their previously reported values. ---
$ws.Range("B7").Value = "..."
$ws.Range("C7").Value = "..."
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 11
$ws.Range("G7").Value = 9
$ws.Range("H7").Value = "..."
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = "..."
$ws.Range("K7").Value = "..."
$ws.Range("L7").Value = "..."
$ws.Range("M7").Value = "..."
$ws.Range("N7").Value = "..."
$ws.Range("O7").Value = 5

# --- New row 8: footnote explaining the "..." marker ---
# A9 already carries the small 9pt Arial "note" font with no fill/border,
# so reuse that formatting for the new footnote cell.
$ws.Range("A9").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$note = "Note: „ ... „ - Data is confidential or unavailable."
$ws.Range("A8").Value = $note
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 9

$bold = $ws.Range("A8").Characters(1, 5)
$bold.Font.Name = "Arial"
$bold.Font.Size = 9
$bold.Font.Bold = $true
$bold.Font.Underline = $true

$rest = $ws.Range("A8").Characters(6, $note.Length - 5)
$rest.Font.Name = "Arial"
$rest.Font.Size = 9

Write-Output "done"
